$d = $word.ActiveDocument

# Paragraph 3 ("Some people use single underlines for emphasis.") gets its
# inline formatting reworked:
#   - "single underlines for" -> "underlining" (underline now a bare <w:u/>
#     covering only that one word, not the trailing space)
#   - the italic+underlined "emphasis" run loses all its direct formatting
#     and is merged with the following "." into a single plain run, with
#     "for " reattached in front of it.
$target = $d.Paragraphs.Item(3).Range

$newXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Some people use</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:u/>
              </w:rPr>
              <w:t xml:space="preserve">underlining</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">for emphasis.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$target.InsertXML($newXml)
